$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '30.687.38'
$ws.Cells.Item(2, 5).Value = '  +0.64%  '

Set-TextValue 3 4 '1.945.89'
$ws.Cells.Item(3, 5).Value = '  +1.64%  '

Set-TextValue 4 4 '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

Set-TextValue 5 4 '247.07'
$ws.Cells.Item(5, 5).Value = '  +0.96%  '

Set-TextValue 6 4 '1.001'

Set-TextValue 7 4 '0.4822'
$ws.Cells.Item(7, 5).Value = '  -0.24%  '

$ws.Cells.Item(8, 2).Value = 'Cardano'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 8 4 '0.2929'
$ws.Cells.Item(8, 5).Value = '  +1.14%  '

$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 9 4 '0.06807'
$ws.Cells.Item(9, 5).Value = '  +1.25%  '

$ws.Cells.Item(10, 2).Value = 'Litecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 10 4 '112.13'
$ws.Cells.Item(10, 5).Value = '  +2.16%  '

$ws.Cells.Item(11, 2).Value = 'Solana'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 11 4 '19.37'
$ws.Cells.Item(11, 5).Value = '  +2.05%  '

$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 12 4 '1.946.05'
$ws.Cells.Item(12, 5).Value = '  +1.65%  '

$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 13 4 '0.07670'
$ws.Cells.Item(13, 5).Value = '  +1.58%  '

$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 14 4 '5.499'
$ws.Cells.Item(14, 5).Value = '  +4.27%  '

$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 15 4 '0.6870'
$ws.Cells.Item(15, 5).Value = '  +2.16%  '

$ws.Cells.Item(16, 2).Value = 'BitcoinCash'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 16 4 '294.99'
$ws.Cells.Item(16, 5).Value = '  +6.28%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 17 4 '30.713.32'
$ws.Cells.Item(17, 5).Value = '  +0.71%  '

$ws.Cells.Item(18, 2).Value = 'Avalanche'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 18 4 '13.24'
$ws.Cells.Item(18, 5).Value = '  +2.95%  '

$ws.Cells.Item(19, 2).Value = 'Uniswap'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 19 4 '5.660'
$ws.Cells.Item(19, 5).Value = '  +2.83%  '

$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 20 4 '0.000007677'
$ws.Cells.Item(20, 5).Value = '  +1.41%  '

$ws.Cells.Item(21, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 21 4 '2.201.68'
$ws.Cells.Item(21, 5).Value = '  +1.63%  '

$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 22 4 '1.001'
$ws.Cells.Item(22, 5).Value = '  +0.08%  '

$ws.Cells.Item(23, 2).Value = 'BinanceUSD'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 23 4 '1.001'
$ws.Cells.Item(23, 5).Value = '  +0.09%  '

$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 24 4 '6.579'
$ws.Cells.Item(24, 5).Value = '  +1.62%  '

$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 25 4 '9.822'
$ws.Cells.Item(25, 5).Value = '  +3.96%  '

$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 26 4 '168.47'
$ws.Cells.Item(26, 5).Value = '  +2.84%  '

$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 27 4 '20.25'
$ws.Cells.Item(27, 5).Value = '  +0.27%  '

$ws.Cells.Item(28, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 28 4 '2.179'
$ws.Cells.Item(28, 5).Value = '  +2.53%  '

$ws.Cells.Item(29, 2).Value = 'Stellar'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 29 4 '0.1079'
$ws.Cells.Item(29, 5).Value = '  +2.29%  '

$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 30 4 '1.436'
$ws.Cells.Item(30, 5).Value = '  +2.51%  '

$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 31 4 '4.690'
$ws.Cells.Item(31, 5).Value = '  +15.67%  '

$ws.Cells.Item(32, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 32 4 '4.501'
$ws.Cells.Item(32, 5).Value = '  +8.26%  '

$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 33 4 '0.05063'
$ws.Cells.Item(33, 5).Value = '  +1.42%  '

$ws.Cells.Item(34, 2).Value = 'ImmutableX'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 34 4 '0.7722'
$ws.Cells.Item(34, 5).Value = '  +5.73%  '

$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 35 4 '1.158'
$ws.Cells.Item(35, 5).Value = '  +2.11%  '

$ws.Cells.Item(36, 2).Value = 'VeChain'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 36 4 '0.02071'
$ws.Cells.Item(36, 5).Value = '  +2.01%  '

$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 37 4 '2.732'
$ws.Cells.Item(37, 5).Value = '  +0.02%  '

$ws.Cells.Item(38, 2).Value = 'MXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 38 4 '2.702'
$ws.Cells.Item(38, 5).Value = '  +1.50%  '

$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 39 4 '2.050'
$ws.Cells.Item(39, 5).Value = '  +1.73%  '

$ws.Cells.Item(40, 2).Value = 'Quant'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 40 4 '111.00'
$ws.Cells.Item(40, 5).Value = '  +0.02%  '

$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 41 4 '0.4452'
$ws.Cells.Item(41, 5).Value = '  +0.23%  '

$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 42 4 '0.8744'
$ws.Cells.Item(42, 5).Value = '  +1.19%  '

$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 43 4 '5.906'
$ws.Cells.Item(43, 5).Value = '  +1.63%  '

$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 44 4 '1.003'
$ws.Cells.Item(44, 5).Value = '  +0.27%  '

Set-TextValue 45 4 '69.67'
$ws.Cells.Item(45, 5).Value = '  +2.53%  '

$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 46 4 '7.352'
$ws.Cells.Item(46, 5).Value = '  -0.03%  '

$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 47 4 '9.354'
$ws.Cells.Item(47, 5).Value = '  +1.35%  '

$ws.Cells.Item(48, 2).Value = 'BitcoinSV'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 48 4 '48.40'
$ws.Cells.Item(48, 5).Value = '  +0.76%  '

$ws.Cells.Item(49, 2).Value = 'Algorand'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 49 4 '0.1252'
$ws.Cells.Item(49, 5).Value = '  +0.88%  '

$ws.Cells.Item(50, 2).Value = 'Elrond'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 50 4 '35.68'
$ws.Cells.Item(50, 5).Value = '  +2.75%  '

$ws.Cells.Item(51, 2).Value = 'WOONetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
Set-TextValue 51 4 '0.2517'
$ws.Cells.Item(51, 5).Value = '  +0.82%  '
